$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the existing data rows (row 2 and row 3) down
# onto the two new rows (4 and 5) before writing values, so the new cells
# pick up the same (yellow-highlighted) styles as the rest of the table.
$ws.Range("A2:R2").Copy() | Out-Null
$ws.Range("A4:R4").PasteSpecial(-4122) | Out-Null
$ws.Range("A3:R3").Copy() | Out-Null
$ws.Range("A5:R5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 4: deleted "species1" animal record, now appears with disposition info
$ws.Range("A4").Value = "Australia"
$ws.Range("B4").Value = 1111113
$ws.Range("C4").Value = "big one"
$ws.Range("D4").Value = "class1"
$ws.Range("E4").Value = "order1"
$ws.Range("F4").Value = "family1"
$ws.Range("G4").Value = "genus1"
$ws.Range("H4").Value = "species1"
$ws.Range("N4").Value = "F"
$ws.Range("O4").Value = "a square nose"
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 0

# Row 5: deleted "species2" animal record, now appears with disposition info
$ws.Range("A5").Value = "Australia"
$ws.Range("B5").Value = 1111114
$ws.Range("C5").Value = "big bubba"
$ws.Range("D5").Value = "class1"
$ws.Range("E5").Value = "order1"
$ws.Range("F5").Value = "family1"
$ws.Range("G5").Value = "genus1"
$ws.Range("H5").Value = "species2"
$ws.Range("N5").Value = "F"
$ws.Range("O5").Value = "a yellow nose"
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 0

$ws.Range("A4").Select()
